$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 12-15 down to 13-16 (bottom-up to avoid clobbering) ---
# Row 15 (Immutable / BooleanType[null]) -> Row 16
$ws.Range("A15:B15").Copy($ws.Range("A16:B16"))
# Row 14 (Copyright / blank) -> Row 15
$ws.Range("A14:B14").Copy($ws.Range("A15:B15"))
$ws.Range("B15").ClearContents()
# Row 13 (Purpose / blank) -> Row 14
$ws.Range("A13:B13").Copy($ws.Range("A14:B14"))
$ws.Range("B14").ClearContents()
# Row 12 (Description / Permissible values...) -> Row 13
$ws.Range("A12:B12").Copy($ws.Range("A13:B13"))

# --- Row 12 becomes the new "Jurisdiction" row ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# --- Row 11 keeps "Contact" label, new value for second contact line ---
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Row 10 Contact value updated ---
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Top metadata field updates ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"
